$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Manthan"
$ws.Range("A10").Value = "Mali"
$ws.Range("A11").Value = 416410

$ws.Range("A12").Select()
